$wb = $excel.ActiveWorkbook

# --- Locate the template sheet (jour5) that the new "jour 6" sheet is based on ---
$src = $wb.Worksheets.Item("jour5")

# Select the full used range before copying, mirroring how the source workbook's
# sheet tab ended up with its selection set to A1:E65 after the duplication.
$src.Activate()
$src.Range("A1:E65").Select()

# Duplicate "jour5" and place the copy right after it (i.e. before "Feuil1")
$src.Copy($null, $src)

# The new sheet is inserted immediately after "jour5" and is named "jour5 (2)"
$newSheet = $wb.Worksheets.Item("jour5 (2)")
$newSheet.Name = "jour 6"

# --- Update the header row for the new scrum day ---
$newSheet.Range("C3").Value = "Thomas Garneau"
$newSheet.Range("D3").Value = "Date : mercredi 2 novembre 2022 "

# --- Paul Agudze ---
$newSheet.Range("C7").Value = "Le timer du jeu, assemblage des derniers détails du menu"
$newSheet.Range("C8").Value = "Lier le canvas de l'input name"
$newSheet.Range("C9").Value = "Non"

# --- Maxime Desrochers ---
$newSheet.Range("C12").Value = ""
$newSheet.Range("C13").Value = "Rien"
$newSheet.Range("C14").Value = "Non"

# --- Thomas Garneau ---
$newSheet.Range("C17").Value = "Arranger le bug de vitesse et documenter code "
$newSheet.Range("C18").Value = "Rien"
$newSheet.Range("C19").Value = "Non"

# --- Nancy Wlodarski ---
$newSheet.Range("C22").Value = "Afficher la fenêtre GameOver "
$newSheet.Range("C23").Value = "Rien"
$newSheet.Range("C24").Value = "Non"

# --- Re-create the "Solutions-Suggestions" table on the new sheet ---
$tbl = $newSheet.ListObjects.Add(1, $newSheet.Range("B42:E65"), $null, 1)
$tbl.Name = "Tableau1468"
$tbl.TableStyle = "TableStyleLight9"

# Put the active cell/selection on the new sheet where the author left it
$newSheet.Activate()
$newSheet.Range("C3").Select()
